# Nouvelle affectation sur la liste des taches
# Affectation de Erwann a :
#   - annulerReservation            (row 22)
#   - CtrlAnnulerReservation        (row 47)
#   - VueAnnulerReservation         (row 48)
# Ces taches passent de "A faire" a "En cours", responsable = Erwann,
# date de debut = 04/10/2016.
#
# Dans la foulee, les taches deja affectees a Erwann et qui sont
# terminees passent de "En cours" a "Termine" avec une date de fin :
#   - getLesSalles                  (row 33)
#   - CtrlConsulterSalles           (row 62)
#   - VueConsulterSalles            (row 63)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dateDebut = Get-Date -Year 2016 -Month 10 -Day 4 -Hour 0 -Minute 0 -Second 0

# --- Row 22 : annulerReservation -> affectee a Erwann, En cours ---
$ws.Range("B22").Value = ""
$ws.Range("C22").Value = "X"
$ws.Range("E22").Value = "Erwann"
$ws.Range("F22").Value = $dateDebut

# --- Row 47 : CtrlAnnulerReservation -> affectee a Erwann, En cours ---
$ws.Range("B47").Value = ""
$ws.Range("C47").Value = "X"
$ws.Range("E47").Value = "Erwann"
$ws.Range("F47").Value = $dateDebut

# --- Row 48 : VueAnnulerReservation -> affectee a Erwann, En cours ---
$ws.Range("B48").Value = ""
$ws.Range("C48").Value = "X"
$ws.Range("E48").Value = "Erwann"
$ws.Range("F48").Value = $dateDebut

# --- Row 33 : getLesSalles -> Termine (deja assignee a Erwann) ---
$ws.Range("C33").Value = ""
$ws.Range("D33").Value = "X"
$ws.Range("F33").Value = $dateDebut
$ws.Range("G33").Value = $dateDebut

# --- Row 62 : CtrlConsulterSalles -> Termine (deja assignee a Erwann) ---
$ws.Range("C62").Value = ""
$ws.Range("D62").Value = "X"
$ws.Range("F62").Value = $dateDebut
$ws.Range("G62").Value = $dateDebut

# --- Row 63 : VueConsulterSalles -> Termine (deja assignee a Erwann) ---
$ws.Range("C63").Value = ""
$ws.Range("D63").Value = "X"
$ws.Range("F63").Value = $dateDebut
$ws.Range("G63").Value = $dateDebut

# --- Vue : selection + scroll sur la zone modifiee ---
$ws.Activate()
$ws.Range("G22").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
